$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.449.97"
$ws.Range("E2").Value = "'  -2.52%  "
$ws.Range("D3").Value = "'2.417.37"
$ws.Range("E3").Value = "'  -4.62%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'572.55"
$ws.Range("E5").Value = "'  -3.69%  "
$ws.Range("D6").Value = "'164.48"
$ws.Range("E6").Value = "'  -7.25%  "
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("D8").Value = "'0.500"
$ws.Range("E8").Value = "'  -5.98%  "
$ws.Range("D9").Value = "'2.415.57"
$ws.Range("E9").Value = "'  -4.64%  "
$ws.Range("E10").Value = "'  -8.19%  "
$ws.Range("E11").Value = "'  -1.37%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'4.75"
$ws.Range("E12").Value = "'  -7.29%  "
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "'0.323"
$ws.Range("E13").Value = "'  -6.76%  "
$ws.Range("D14").Value = "'24.68"
$ws.Range("E14").Value = "'  -8.08%  "
$ws.Range("D15").Value = "'66.400.35"
$ws.Range("E15").Value = "'  -2.45%  "
$ws.Range("D16").Value = "'2.742.17"
$ws.Range("E16").Value = "'  -8.44%  "
$ws.Range("D17").Value = "'0.0000165"
$ws.Range("E17").Value = "'  -8.57%  "
$ws.Range("D18").Value = "'2.384.34"
$ws.Range("E18").Value = "'  -5.85%  "
$ws.Range("D19").Value = "'11.00"
$ws.Range("E19").Value = "'  -4.95%  "
$ws.Range("E20").Value = "'  -8.53%  "
$ws.Range("D21").Value = "'346.98"
$ws.Range("E21").Value = "'  -5.56%  "
$ws.Range("D22").Value = "'3.96"
$ws.Range("E22").Value = "'  -5.85%  "
$ws.Range("E23").Value = "'  +0.08%  "
$ws.Range("D24").Value = "'68.11"
$ws.Range("E24").Value = "'  -4.11%  "
$ws.Range("D25").Value = "'4.12"
$ws.Range("E25").Value = "'  -12.47%  "
$ws.Range("E26").Value = "'  -10.93%  "
$ws.Range("D27").Value = "'8.87"
$ws.Range("E27").Value = "'  -13.05%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "'  -0.09%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "'2.520.57"
$ws.Range("E29").Value = "'  -5.11%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0₃0870"
$ws.Range("E30").Value = "'  -12.97%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'7.64"
$ws.Range("E31").Value = "'  -7.81%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'485.41"
$ws.Range("E32").Value = "'  -10.53%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.75"
$ws.Range("E33").Value = "'  -7.09%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "'  -0.11%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.20"
$ws.Range("E35").Value = "'  -10.80%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'156.23"
$ws.Range("E36").Value = "'  -0.66%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.112"
$ws.Range("E37").Value = "'  -13.12%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'18.51"
$ws.Range("E38").Value = "'  -1.01%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'18.08"
$ws.Range("E39").Value = "'  -4.14%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'1.32"
$ws.Range("E40").Value = "'  -9.93%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.64"
$ws.Range("E41").Value = "'  -9.55%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.321"
$ws.Range("E42").Value = "'  -9.91%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "'4.57"
$ws.Range("E43").Value = "'  -12.24%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'39.03"
$ws.Range("E44").Value = "'  -2.23%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "'  -10.49%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'137.20"
$ws.Range("E46").Value = "'  -7.00%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'3.42"
$ws.Range("E47").Value = "'  -8.51%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.503"
$ws.Range("E48").Value = "'  -10.26%  "
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "'1.55"
$ws.Range("E49").Value = "'  -9.15%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0717"
$ws.Range("E50").Value = "'  -5.37%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.573"
$ws.Range("E51").Value = "'  -4.15%  "
